# Update countries & provincias Spain
# Reflects a new data pull: "Datos actualizados" timestamp moves from 03:35 to 04:05,
# Bolivia overtakes Luxemburgo (rows 67/68 swap), Curazao overtakes Dominica
# (rows 201/202 swap), and a handful of other countries get refreshed counts
# (Honduras row 77, El Salvador row 98).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Title / timestamp row ---
$ws.Range("A1").Value = "Datos actualizados a 18 de Mayo de 2020 a las 04:05"

# --- Rows 67/68: Bolivia moves above Luxemburgo ---
$ws.Range("A67").Value = "Bolivia"
$ws.Range("B67").Value = 4088
$ws.Range("C67").Value = 262
$ws.Range("D67").Value = 493
$ws.Range("E67").Value = 3426
$ws.Range("F67").Value = 0
$ws.Range("G67").Value = 4
$ws.Range("H67").Value = 169

$ws.Range("A68").Value = "Luxemburgo"
$ws.Range("B68").Value = 3945
$ws.Range("C68").Value = 0
$ws.Range("D68").Value = 3702
$ws.Range("E68").Value = 136
$ws.Range("F68").Value = 0
$ws.Range("G68").Value = 0
$ws.Range("H68").Value = 107

# --- Row 77: Honduras refreshed counts ---
$ws.Range("B77").Value = 2646
$ws.Range("C77").Value = 81
$ws.Range("D77").Value = 319
$ws.Range("E77").Value = 2185
$ws.Range("F77").Value = 0
$ws.Range("G77").Value = 4
$ws.Range("H77").Value = 142

# --- Row 98: El Salvador refreshed counts ---
$ws.Range("D98").Value = 464
$ws.Range("E98").Value = 844
$ws.Range("G98").Value = 3
$ws.Range("H98").Value = 30

# --- Rows 201/202: Curazao moves above Dominica ---
$ws.Range("A201").Value = "Curazao"
$ws.Range("B201").Value = 16
$ws.Range("C201").Value = 0
$ws.Range("D201").Value = 14
$ws.Range("E201").Value = 1
$ws.Range("F201").Value = 0
$ws.Range("G201").Value = 0
$ws.Range("H201").Value = 1

$ws.Range("A202").Value = "Dominica"
$ws.Range("B202").Value = 16
$ws.Range("C202").Value = 0
$ws.Range("D202").Value = 16
$ws.Range("E202").Value = 0
$ws.Range("F202").Value = 0
$ws.Range("G202").Value = 0
$ws.Range("H202").Value = 0
